# Update "想去人数" (column F) values across sheets to match the
# regenerated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rId1 / sheet1.xml) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 2442
$ws1.Range("F8").Value  = 375
$ws1.Range("F9").Value  = 3400
$ws1.Range("F10").Value = 899
$ws1.Range("F11").Value = 604
$ws1.Range("F12").Value = 856
$ws1.Range("F13").Value = 1543
$ws1.Range("F15").Value = 1
$ws1.Range("F16").Value = 969
$ws1.Range("F17").Value = 1745
$ws1.Range("F20").Value = 1521
$ws1.Range("F22").Value = 90
$ws1.Range("F24").Value = 4077
$ws1.Range("F26").Value = 2659
$ws1.Range("F27").Value = 1184

# --- Sheet "演出" (rId2 / sheet2.xml) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F20").Value = 1
$ws2.Range("F27").Value = 56
$ws2.Range("F37").Value = 1
$ws2.Range("F47").Value = 12
$ws2.Range("F48").Value = 12
$ws2.Range("F49").Value = 312

# --- Sheet "本地生活" (rId3 / sheet3.xml) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value  = 2520
$ws3.Range("F7").Value  = 9557
$ws3.Range("F12").Value = 2877
$ws3.Range("F13").Value = 410
$ws3.Range("F14").Value = 734
$ws3.Range("F15").Value = 100

# --- Sheet "全部类型" (rId4 / sheet4.xml) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 410
$ws4.Range("F9").Value  = 734
$ws4.Range("F11").Value = 100
$ws4.Range("F16").Value = 375
$ws4.Range("F17").Value = 3400
$ws4.Range("F19").Value = 899
$ws4.Range("F20").Value = 604
$ws4.Range("F21").Value = 856
$ws4.Range("F25").Value = 969
$ws4.Range("F29").Value = 1745
$ws4.Range("F32").Value = 1521
$ws4.Range("F34").Value = 56
$ws4.Range("F38").Value = 90
$ws4.Range("F41").Value = 4077
$ws4.Range("F43").Value = 2659
$ws4.Range("F48").Value = 12
$ws4.Range("F49").Value = 1184
